# Insert a new weekly price record for Mango (Vega Monumental Concepción) at
# row 95, pushing the existing rows 95-167 down to 96-168.
#
# The new row carries the same fixed/categorical data as the rest of the
# sheet (market, region, product classification, unit, origin country,
# quality grade, kg/unit), but with its own date + price figures, so the
# simplest reliable way to create it is to insert a blank row, clone the
# row immediately below it (which, post-insert, holds the old row 95 data)
# into the new row, and then overwrite just the cells that differ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 95:167 down to 96:168.
$ws.Rows(95).Insert()

# Seed the new row 95 with the same row layout/content as its neighbour
# (now row 96, which holds what used to be row 95), then correct the
# cells that are actually new.
$ws.Range("A96:T96").Copy()
$ws.Range("A95:T95").PasteSpecial()

$ws.Range("D95").Value = 45072
$ws.Range("M95").Value = 200
$ws.Range("N95").Value = 7500
$ws.Range("O95").Value = 8000
$ws.Range("P95").Value = 7750
$ws.Range("S95").Value = 1938
